$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 text (Conversión del día) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.73 = 61654.27 pesos`n✅ 61654.27 pesos = 14.63 = 969.98 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update "tasas" sheet rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 67.879
$wsTasas.Range("O10").Value = 4185.03
$wsTasas.Range("N12").Value = 4215
$wsTasas.Range("O12").Value = 66.313
